$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-18 Thursday" "2024-01-19 Friday"

Replace-Text "443÷8=" "674÷6="
Replace-Text "681÷4=" "946÷9="
Replace-Text "749÷5=" "915÷6="
Replace-Text "494÷6=" "870÷5="
Replace-Text "749÷2=" "179÷4="

Replace-Text "442÷8=" "637÷4="
Replace-Text "947÷5=" "859÷7="
Replace-Text "881÷4=" "489÷8="
Replace-Text "334÷8=" "314÷5="
Replace-Text "886÷9=" "163÷6="

Replace-Text "111÷8=" "628÷2="
Replace-Text "549÷5=" "303÷2="
Replace-Text "577÷2=" "749÷7="
Replace-Text "874÷8=" "986÷8="
Replace-Text "184÷5=" "857÷5="

Replace-Text "797÷3=" "232÷9="
Replace-Text "325÷8=" "862÷4="
Replace-Text "507÷2=" "435÷4="
Replace-Text "150÷7=" "886÷7="
Replace-Text "273÷7=" "903÷5="

Replace-Text "449÷9=" "397÷9="
Replace-Text "430÷2=" "763÷4="
Replace-Text "426÷3=" "304÷4="
Replace-Text "334÷4=" "140÷6="
Replace-Text "837÷5=" "631÷2="
